$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "79.868.19"
$ws.Range("E2").Value = "  +5.11%  "
$ws.Range("D3").Value = "3.216.88"
$ws.Range("E3").Value = "  +6.71%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'211.56"
$ws.Range("E5").Value = "  +7.38%  "
$ws.Range("D6").Value = "'639.61"
$ws.Range("E6").Value = "  +3.89%  "
$ws.Range("D7").Value = "'0.268"
$ws.Range("E7").Value = "  +30.75%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  +10.62%  "
$ws.Range("D10").Value = "3.216.10"
$ws.Range("E10").Value = "  +6.76%  "
$ws.Range("D12").Value = "'0.0000273"
$ws.Range("E12").Value = "  +43.35%  "
$ws.Range("E13").Value = "  +3.73%  "
$ws.Range("E14").Value = "  +4.66%  "
$ws.Range("D15").Value = "3.800.70"
$ws.Range("E15").Value = "  +6.46%  "
$ws.Range("D16").Value = "'32.77"
$ws.Range("E16").Value = "  +13.87%  "
$ws.Range("D17").Value = "79.609.63"
$ws.Range("E17").Value = "  +4.85%  "
$ws.Range("D18").Value = "3.199.04"
$ws.Range("E18").Value = "  +6.34%  "
$ws.Range("B20").Value = "SuiNetwork"
$ws.Range("C20").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D20").Value = "'3.04"
$ws.Range("E20").Value = "  +28.90%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'9.41"
$ws.Range("E21").Value = "  +6.10%  "
$ws.Range("D22").Value = "'447.54"
$ws.Range("E22").Value = "  +18.63%  "
$ws.Range("D23").Value = "'5.35"
$ws.Range("E23").Value = "  +22.70%  "
$ws.Range("E24").Value = "  +13.25%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "'77.96"
$ws.Range("E25").Value = "  +8.05%  "
$ws.Range("B26").Value = "WrappedeETH"
$ws.Range("C26").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D26").Value = "3.366.32"
$ws.Range("E26").Value = "  +6.14%  "
$ws.Range("D27").Value = "'10.97"
$ws.Range("E27").Value = "  +12.77%  "
$ws.Range("E28").Value = "  +0.13%  "
$ws.Range("D29").Value = "'0.0000127"
$ws.Range("E29").Value = "  +18.90%  "
$ws.Range("E30").Value = "  +13.24%  "
$ws.Range("D31").Value = "'0.999"
$ws.Range("E31").Value = "  +0.25%  "
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").Value = "'1.54"
$ws.Range("E32").Value = "  +11.34%  "
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").Value = "'565.36"
$ws.Range("E33").Value = "  +15.03%  "
$ws.Range("D34").Value = "'0.158"
$ws.Range("E34").Value = "  +30.34%  "
$ws.Range("E35").Value = "  +7.69%  "
$ws.Range("D36").Value = "'23.54"
$ws.Range("E36").Value = "  +14.91%  "
$ws.Range("D37").Value = "'0.124"
$ws.Range("E37").Value = "  +20.66%  "
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("D39").Value = "'0.416"
$ws.Range("E39").Value = "  +10.75%  "
$ws.Range("D40").Value = "'163.73"
$ws.Range("E40").Value = "  +1.09%  "
$ws.Range("D41").Value = "'5.83"
$ws.Range("E41").Value = "  +14.60%  "
$ws.Range("D43").Value = "'192.42"
$ws.Range("E43").Value = "  +1.21%  "
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("E45").Value = "  +12.55%  "
$ws.Range("E46").Value = "  +14.28%  "
$ws.Range("D47").Value = "'0.804"
$ws.Range("E47").Value = "  +4.18%  "
$ws.Range("E48").Value = "  +9.01%  "
$ws.Range("E49").Value = "  +13.32%  "
$ws.Range("D50").Value = "'43.14"
$ws.Range("E50").Value = "  +4.94%  "
$ws.Range("B51").Value = "ARBITRUM"
$ws.Range("C51").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D51").Value = "'0.647"
$ws.Range("E51").Value = "  +9.67%  "
